# Update crypto price/volume figures per the Aug 4 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.622.05"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "2.904.11"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'526.89"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").Value = "'143.23"
$ws.Range("E6").Value = "  -5.02%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -3.41%  "
$ws.Range("D9").Value = "2.911.25"
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("E10").Value = "  -4.90%  "
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").Value = "'0.359"
$ws.Range("E12").Value = "  -2.50%  "
$ws.Range("D13").Value = "3.409.50"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").Value = "'0.128"
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").Value = "60.607.64"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").Value = "'22.48"
$ws.Range("E16").Value = "  -4.79%  "
$ws.Range("D17").Value = "2.907.58"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("E18").Value = "  -3.94%  "
$ws.Range("D19").Value = "'4.97"
$ws.Range("E19").Value = "  -3.38%  "
$ws.Range("D20").Value = "'11.59"
$ws.Range("E20").Value = "  -3.84%  "
$ws.Range("D21").Value = "'350.48"
$ws.Range("E21").Value = "  -7.81%  "
$ws.Range("D22").Value = "'6.55"
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("D25").Value = "'64.69"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("E26").Value = "  -4.13%  "
$ws.Range("E27").Value = "  -5.92%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").Value = "'7.80"
$ws.Range("E29").Value = "  -4.39%  "
$ws.Range("D30").Value = "0.0₃0854"
$ws.Range("E30").Value = "  -8.53%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("E33").Value = "  -4.55%  "
$ws.Range("D34").Value = "'152.16"
$ws.Range("E34").Value = "  -4.18%  "
$ws.Range("D35").Value = "'4.30"
$ws.Range("E35").Value = "  -5.79%  "
$ws.Range("E36").Value = "  -5.80%  "
$ws.Range("D37").Value = "'0.995"
$ws.Range("E37").Value = "  -6.46%  "
$ws.Range("E38").Value = "  -5.86%  "
$ws.Range("D39").Value = "'37.61"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").Value = "'1.47"
$ws.Range("E40").Value = "  -4.59%  "
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("D42").Value = "2.289.62"
$ws.Range("E42").Value = "  -5.01%  "
$ws.Range("E43").Value = "  -3.23%  "
$ws.Range("D44").Value = "'0.0579"
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("D45").Value = "'20.30"
$ws.Range("E45").Value = "  -7.86%  "
$ws.Range("D46").Value = "'0.997"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'4.94"
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("E50").Value = "  -3.85%  "
$ws.Range("D51").Value = "'18.26"
$ws.Range("E51").Value = "  -7.37%  "
